$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update forecast error table values for rows 7-11 (Q5-Q9)
$ws.Range("B7").Value = 0.04537116952229212
$ws.Range("C7").Value = 0.9472648544446646
$ws.Range("D7").Value = 2.781695858718539
$ws.Range("E7").Value = 1.667841676754283
$ws.Range("F7").Value = 1.689604287563018
$ws.Range("G7").Value = 38

$ws.Range("B8").Value = 0.1157334130921949
$ws.Range("C8").Value = 1.016311866239339
$ws.Range("D8").Value = 2.967633790164031
$ws.Range("E8").Value = 1.722682150068326
$ws.Range("F8").Value = 1.742498716630711
$ws.Range("G8").Value = 37

$ws.Range("B9").Value = -0.02432183078563367
$ws.Range("C9").Value = 1.165694726323169
$ws.Range("D9").Value = 4.494950312546154
$ws.Range("E9").Value = 2.120129786722066
$ws.Range("F9").Value = 2.17506412757088
$ws.Range("G9").Value = 20

$ws.Range("B10").Value = -0.4137243084650019
$ws.Range("C10").Value = 0.9127267147502702
$ws.Range("D10").Value = 1.863584276199294
$ws.Range("E10").Value = 1.365131596659932
$ws.Range("F10").Value = 1.354049917414355
$ws.Range("G10").Value = 13

$ws.Range("B11").Value = 0.1972777442226493
$ws.Range("C11").Value = 0.547588270575837
$ws.Range("D11").Value = 0.4683298695408837
$ws.Range("E11").Value = 0.6843463081955536
$ws.Range("F11").Value = 0.7326419326445446
$ws.Range("G11").Value = 5
